# Update the cryptocurrency list: refresh Price (column D) and Volume(1h)
# percentage-change (column E) values per the Dec 23 2023 data refresh.
# A leading apostrophe is used for purely numeric-looking Price values so
# Excel keeps them stored as text (matching the original text cells)
# instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.777.72"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.287.86"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'96.32"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").Value = "'269.47"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").Value = "'0.618"
$ws.Range("E7").Value = "  -1.27%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").Value = "'45.21"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("D11").Value = "'0.0934"
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").Value = "'7.89"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").Value = "'15.75"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").Value = "2.631.34"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").Value = "'0.851"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "2.285.05"
$ws.Range("E17").Value = "  -2.16%  "
$ws.Range("D18").Value = "43.753.01"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("E19").Value = "  +3.60%  "
$ws.Range("D20").Value = "'6.18"
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("D21").Value = "'72.10"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("E22").Value = "  +10.82%  "
$ws.Range("D23").Value = "'232.67"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("E24").Value = "  -4.59%  "
$ws.Range("E25").Value = "  +7.15%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'11.37"
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "'38.61"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").Value = "'175.35"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("D32").Value = "'21.81"
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "'5.42"
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("D36").Value = "'4.70"
$ws.Range("E36").Value = "  +8.26%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("D39").Value = "'3.59"
$ws.Range("E39").Value = "  +6.13%  "
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("E41").Value = "  -2.44%  "
$ws.Range("D42").Value = "'12.25"
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("D44").Value = "'64.31"
$ws.Range("E44").Value = "  +4.44%  "
$ws.Range("E45").Value = "  -3.06%  "
$ws.Range("D46").Value = "'0.102"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").Value = "'8.69"
$ws.Range("E47").Value = "  -4.33%  "
$ws.Range("D48").Value = "'97.96"
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("D50").Value = "'0.441"
$ws.Range("E50").Value = "  +6.89%  "
$ws.Range("E51").Value = "  +10.48%  "

